# Updates the "Price" (D) and "Volume(1h)" (E) columns of the cryptos
# worksheet with freshly scraped values. Every touched cell stores a plain
# text string (prices/percentages are kept as text in the source data, e.g.
# "51.589.83" or "  +2.09%  "), so NumberFormat is forced to Text ("@")
# before the write to stop Excel from reinterpreting number/date-looking
# strings, and the explicit style is cleared back to "Normal" afterwards so
# the cell keeps its original (default) styling.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = "51.589.83" },
    @{ Cell = "D3"; Value = "3.018.33" },
    @{ Cell = "E3"; Value = "  +2.09%  " },
    @{ Cell = "D4"; Value = "0.999" },
    @{ Cell = "E4"; Value = "  -0.05%  " },
    @{ Cell = "D5"; Value = "378.72" },
    @{ Cell = "E5"; Value = "  -0.14%  " },
    @{ Cell = "D6"; Value = "102.37" },
    @{ Cell = "E6"; Value = "  -0.09%  " },
    @{ Cell = "D7"; Value = "0.546" },
    @{ Cell = "E7"; Value = "  +0.31%  " },
    @{ Cell = "E9"; Value = "  +0.59%  " },
    @{ Cell = "D10"; Value = "36.70" },
    @{ Cell = "E10"; Value = "  +0.80%  " },
    @{ Cell = "E11"; Value = "  -0.15%  " },
    @{ Cell = "E12"; Value = "  +1.23%  " },
    @{ Cell = "D13"; Value = "3.497.86" },
    @{ Cell = "E13"; Value = "  +1.99%  " },
    @{ Cell = "D14"; Value = "18.41" },
    @{ Cell = "E14"; Value = "  -0.05%  " },
    @{ Cell = "D15"; Value = "7.71" },
    @{ Cell = "E15"; Value = "  -0.33%  " },
    @{ Cell = "D16"; Value = "3.022.17" },
    @{ Cell = "E16"; Value = "  +2.39%  " },
    @{ Cell = "D17"; Value = "0.972" },
    @{ Cell = "E17"; Value = "  -3.94%  " },
    @{ Cell = "D18"; Value = "10.58" },
    @{ Cell = "E18"; Value = "  -14.86%  " },
    @{ Cell = "D19"; Value = "51.565.20" },
    @{ Cell = "E19"; Value = "  +0.96%  " },
    @{ Cell = "D20"; Value = "3.09" },
    @{ Cell = "E20"; Value = "  +0.52%  " },
    @{ Cell = "D21"; Value = "12.43" },
    @{ Cell = "E21"; Value = "  +0.17%  " },
    @{ Cell = "D22"; Value = "0.0₃0963" },
    @{ Cell = "E22"; Value = "  +0.86%  " },
    @{ Cell = "D23"; Value = "69.93" },
    @{ Cell = "E23"; Value = "  +0.37%  " },
    @{ Cell = "D24"; Value = "267.18" },
    @{ Cell = "E24"; Value = "  -0.10%  " },
    @{ Cell = "D25"; Value = "3.15" },
    @{ Cell = "E25"; Value = "  -5.95%  " },
    @{ Cell = "D26"; Value = "8.30" },
    @{ Cell = "E26"; Value = "  +4.12%  " },
    @{ Cell = "D27"; Value = "7.56" },
    @{ Cell = "E27"; Value = "  +8.21%  " },
    @{ Cell = "E28"; Value = "  +3.95%  " },
    @{ Cell = "E29"; Value = "  +0.03%  " },
    @{ Cell = "D30"; Value = "26.16" },
    @{ Cell = "E30"; Value = "  +1.42%  " },
    @{ Cell = "E31"; Value = "  +0.30%  " },
    @{ Cell = "D32"; Value = "10.24" },
    @{ Cell = "E32"; Value = "  -2.71%  " },
    @{ Cell = "D33"; Value = "2.11" },
    @{ Cell = "E33"; Value = "  +3.49%  " },
    @{ Cell = "D34"; Value = "50.55" },
    @{ Cell = "E34"; Value = "  -0.46%  " },
    @{ Cell = "D35"; Value = "33.81" },
    @{ Cell = "E35"; Value = "  -0.68%  " },
    @{ Cell = "D36"; Value = "0.0449" },
    @{ Cell = "E36"; Value = "  +3.33%  " },
    @{ Cell = "E37"; Value = "  -0.12%  " },
    @{ Cell = "E38"; Value = "  +2.57%  " },
    @{ Cell = "D39"; Value = "0.291" },
    @{ Cell = "E39"; Value = "  +13.20%  " },
    @{ Cell = "E40"; Value = "  +0.77%  " },
    @{ Cell = "E41"; Value = "  +1.34%  " },
    @{ Cell = "E42"; Value = "  -0.47%  " },
    @{ Cell = "D43"; Value = "127.25" },
    @{ Cell = "E43"; Value = "  +7.32%  " },
    @{ Cell = "E44"; Value = "  +1.70%  " },
    @{ Cell = "D45"; Value = "3.78" },
    @{ Cell = "E45"; Value = "  +5.82%  " },
    @{ Cell = "D46"; Value = "21.53" },
    @{ Cell = "E46"; Value = "  -0.23%  " },
    @{ Cell = "E48"; Value = "  +2.34%  " },
    @{ Cell = "D49"; Value = "2.025.16" },
    @{ Cell = "E49"; Value = "  -0.76%  " },
    @{ Cell = "D50"; Value = "3.316.37" },
    @{ Cell = "E50"; Value = "  +2.04%  " },
    @{ Cell = "D51"; Value = "0.0317" },
    @{ Cell = "E51"; Value = "  -1.47%  " }
)

foreach ($update in $updates) {
    $range = $ws.Range($update.Cell)
    $range.NumberFormat = "@"
    $range.Value = $update.Value
    $range.Style = "Normal"
}
